$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Capture existing text/number blocks that move to new rows, before
# we overwrite their old locations.
# ------------------------------------------------------------------
$implLabel = $ws.Range("B15").Value2
$impl1Num  = $ws.Range("B16").Value2
$impl1Txt  = $ws.Range("C16").Value2
$impl2Num  = $ws.Range("B17").Value2
$impl2Txt  = $ws.Range("C17").Value2
$impl3Num  = $ws.Range("B18").Value2
$impl4Num  = $ws.Range("B19").Value2

$testsLabel = $ws.Range("B22").Value2
$test1Num   = $ws.Range("B23").Value2
$test1Txt   = $ws.Range("C23").Value2
$test2Num   = $ws.Range("B24").Value2
$test2Txt   = $ws.Range("C24").Value2

$dynNote = $ws.Range("B10").Value2

# Clear out the old locations that will no longer hold data in the new layout
$ws.Range("B10").ClearContents()
$ws.Range("B15:C19").ClearContents()
$ws.Range("B22:C24").ClearContents()

# ------------------------------------------------------------------
# Re-create the moved blocks in their new rows (numbers first so we
# don't disturb shared-string ordering yet).
# ------------------------------------------------------------------
$ws.Range("B33").Value = $testsLabel
$ws.Range("B34").Value = $test1Num
$ws.Range("C34").Value = $test1Txt
$ws.Range("B35").Value = $test2Num
$ws.Range("C35").Value = $test2Txt

$ws.Range("B21").Value = $implLabel
$ws.Range("B22").Value = $impl1Num
$ws.Range("C22").Value = $impl1Txt
$ws.Range("B23").Value = $impl2Num
$ws.Range("C23").Value = $impl2Txt
$ws.Range("B24").Value = $impl3Num
$ws.Range("B25").Value = $impl4Num

# ------------------------------------------------------------------
# Newly authored strings - entered in the same order the original
# author typed them in, so the shared-string table lines up.
# ------------------------------------------------------------------
$ws.Range("C24").Value = "Frequency-Optimized, otherwise matching implementation #1"
$ws.Range("C30").Value = "Maximum Frequency - No cache, 5 stage pipeline, no mulplier/divider/caches. Intended to reach the maximum frequency"
$ws.Range("C27").Value = "Minimum Area, no multiplier, divider, or branch predictor (3 Stage pipeline)"
$ws.Range("C25").Value = "Maximum Frequency, otherwise matching implementation #1"
$ws.Range("F8").Value  = "**"
$ws.Range("B20").Value = "** Dynamic power does not include power from PLL, which makes up ~101mW"
$ws.Range("C28").Value = "Frequency-Optimized"
$ws.Range("C29").Value = "Maximum Performance - No i/d cache, 5 stage pipeline, multiper + branch cache "
$ws.Range("C31").Value = "Area Optimized, includes mulitplier and barrel shifter (3 Stage pipeline)"
$ws.Range("A8").Value  = "2nd tests"
$ws.Range("B26").Value = "pt2"

# ------------------------------------------------------------------
# Remaining numeric fill-in for the "pt2" implementation block.
# ------------------------------------------------------------------
$ws.Range("B27").Value = 1
$ws.Range("B28").Value = 2
$ws.Range("B29").Value = 3
$ws.Range("B30").Value = 4
$ws.Range("B31").Value = 5

# --- Move the "*Dynamic power..." note back in at its new row (19) ---
$ws.Range("B19").Value = $dynNote

# ------------------------------------------------------------------
# First results table (rows 4-7): new "Max Freq" column + row 6 data.
# ------------------------------------------------------------------
$ws.Range("G4").Value = 100
$ws.Range("G5").Value = 100

$ws.Range("B6").Value = 2365
$ws.Range("C6").Value = 2420
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 75
$ws.Range("F6").Value = 21
$ws.Range("G6").Value = 100
$ws.Range("H6").Value = 1.22
$ws.Range("I6").Value = 1.296

# J4 gets its own (non-shared) formula; J5:J13 become one shared formula.
$ws.Range("J4").Formula = "=H4+I4"
$ws.Range("J5:J13").Formula = "=H5+I5"

# ------------------------------------------------------------------
# New second results table (rows 9-13, "2nd tests").
# ------------------------------------------------------------------
$ws.Range("A9").Value = 1
$ws.Range("B9").Value = 1530
$ws.Range("C9").Value = 1359
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 75
$ws.Range("F9").Value = 49
$ws.Range("G9").Value = 125
$ws.Range("H9").Value = 0.446
$ws.Range("I9").Value = 0.565

$ws.Range("A10").Value = 2
$ws.Range("B10").Value = 3166
$ws.Range("C10").Value = 3222
$ws.Range("D10").Value = 3
$ws.Range("E10").Value = 75
$ws.Range("F10").Value = 57
$ws.Range("G10").Value = 125
$ws.Range("H10").Value = 0.655
$ws.Range("I10").Value = 0.832

$ws.Range("A11").Value = 3
$ws.Range("B11").Value = 2316
$ws.Range("C11").Value = 2071
$ws.Range("D11").Value = 3
$ws.Range("E11").Value = 75
$ws.Range("F11").Value = 43
$ws.Range("G11").Value = 100
$ws.Range("H11").Value = 0.525
$ws.Range("I11").Value = 0.663
$ws.Range("J11").Formula = "=H11+I11"

$ws.Range("A12").Value = 4
$ws.Range("B12").Value = 1881
$ws.Range("C12").Value = 1717
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 75
$ws.Range("F12").Value = 69
$ws.Range("G12").Value = 125
$ws.Range("H12").Value = 0.483
$ws.Range("I12").Value = 0.612

$ws.Range("A13").Value = 5
$ws.Range("B13").Value = 1739
$ws.Range("C13").Value = 1510
$ws.Range("D13").Value = 3
$ws.Range("E13").Value = 75
$ws.Range("F13").Value = 68
$ws.Range("G13").Value = 125
$ws.Range("H13").Value = 0.444
$ws.Range("I13").Value = 0.563

# ------------------------------------------------------------------
# Styling: strikethrough font for the first table and the relocated
# implementation block / note that sit alongside it visually.
# ------------------------------------------------------------------
$ws.Range("A4:J7").Font.Strikethrough = $true
$ws.Range("B19").Font.Strikethrough = $true
$ws.Range("B22:C25").Font.Strikethrough = $true

# ------------------------------------------------------------------
# Sheet view changes: zoom + new active selection.
# ------------------------------------------------------------------
$ws.Application.ActiveWindow.Zoom = 160
$ws.Range("D15").Select()

# --- Page setup: portrait orientation ---
$ws.PageSetup.Orientation = 1
